$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 189; existing rows 189-199 shift down to 190-200.
$ws.Rows.Item(189).Insert()

# Populate the new row 189 with data (copying the constant columns from
# the surrounding block, and setting the row-specific values).
$ws.Range("A189").Value = 4
$ws.Range("B189").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C189").Value = "Los Lagos"
$ws.Range("D189").Value = 44509
$ws.Range("E189").Value = 10
$ws.Range("F189").Value = 100112045
$ws.Range("G189").Value = "Zapallo"
$ws.Range("H189").Value = "Paine"
$ws.Range("I189").Value = "1a (guarda)"
$ws.Range("J189").Value = 1100
$ws.Range("K189").Value = 400
$ws.Range("L189").Value = 400
$ws.Range("M189").Value = 400
$ws.Range("N189").Value = "$/kilo (volumen en unidades)"
$ws.Range("O189").Value = "Región de O'Higgins"
$ws.Range("P189").Value = 400
$ws.Range("Q189").Value = 1
$ws.Range("R189").Value = "Hortaliza"

$ws.Range("D189").Style = $ws.Range("D190").Style
